$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "João Rodrigues-Desenho Técnico"
$ws.Range("F3").Value = "Andre Lucca-Circuitos Elétricos"
$ws.Range("E4").Value = "José Ferreira-Tecnologia dos Materiais"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("F6").Value = "-"
